# Edit teste12.docx: capitalize the leading "t" of the existing word and
# append "   minha picccccccccccccccccc       aaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaa",
# reproducing the exact run/proofErr layout Word's spell-checker leaves behind
# (two "misspelled word" wrappers around "Tsteofosf...jj" and "picccc...c").

$d = $word.ActiveDocument

# The document is a single paragraph; grab it (and remember its own
# paragraph-mark identity so the rebuilt paragraph keeps the same
# w14:paraId / w14:textId / rsid attributes it already had).
$para = $d.Paragraphs(1)
$range = $para.Range

$wNs   = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$w14Ns = "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

$newParagraphXml = "<w:p $wNs $w14Ns w14:paraId='388AF7C5' w14:textId='490E52C6' w:rsidR='008E493B' w:rsidRDefault='00B9718D'>" +
  "<w:proofErr w:type='spellStart'/>" +
  "<w:r><w:t>T</w:t></w:r>" +
  "<w:r><w:t>steofosfjjjjjjjjjjjjjjjjjjjjjjjjjjjjjjjjjjjj</w:t></w:r>" +
  "<w:proofErr w:type='spellEnd'/>" +
  "<w:r><w:t xml:space='preserve'>   minha </w:t></w:r>" +
  "<w:proofErr w:type='spellStart'/>" +
  "<w:r><w:t>picccccccccccccccccc</w:t></w:r>" +
  "<w:proofErr w:type='spellEnd'/>" +
  "<w:r><w:t xml:space='preserve'>       aaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaa</w:t></w:r>" +
  "</w:p>"

# InsertXML replaces the contents of the exact range it targets; targeting
# the whole paragraph range swaps the single run for the run/proofErr
# sequence Word leaves after the user fixes the capitalisation and keeps
# typing past it.
$result = $range.InsertXML($newParagraphXml)
